$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data as scraped on Fri Aug  9 05:30:55 UTC 2024

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.878.55"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +6.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.665.09"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +9.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.663.56"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +9.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.37"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +11.28%  "
$ws.Range("E11").Value = "  +5.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("D12").ClearFormats()
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.115.08"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +10.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.919.66"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.78"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.99%  "
$ws.Range("E17").Value = "  +4.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.656.57"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +9.63%  "
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +7.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.48"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.35%  "
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.39%  "
$ws.Range("E25").Value = "  +3.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.752.41"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +9.37%  "
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0859"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +9.57%  "
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.48"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.91%  "
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("E35").Value = "  +6.43%  "
$ws.Range("E36").Value = "  +8.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.21"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.52"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.858"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "307.79"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +14.10%  "
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.837"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +27.28%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.76"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.34"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.642"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0579"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.96%  "
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.90"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +14.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.88"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0236"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.035.22"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.10%  "
